# after running cases A-F for the first time
# Updates the numeric results in row 2 of sheet "A-03" with the freshly
# computed values from the simulation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6993.552241803531
$ws.Range("C2").Value = 13330.70862908129
$ws.Range("D2").Value = -5686.582311710362
$ws.Range("E2").Value = -650.5740755674435
$ws.Range("F2").Value = 60.92981298656036
$ws.Range("G2").Value = 45.0390723501958
$ws.Range("H2").Value = 61.37293967784098
$ws.Range("I2").Value = 44.69642864647369
$ws.Range("J2").Value = 60.83066004025909
$ws.Range("K2").Value = 45.00241456369662
$ws.Range("L2").Value = 54.99078891686219
$ws.Range("O2").Value = 45.95613052988386
$ws.Range("P2").Value = 53.12368901672909
$ws.Range("R2").Value = 6.287176378676221
$ws.Range("S2").Value = -11.41957657347872
$ws.Range("T2").Value = 5.132400194802501
$ws.Range("X2").Value = -45.53134160081841
$ws.Range("Y2").Value = -70.18623930993614
$ws.Range("Z2").Value = -53.75641409691111
$ws.Range("AE2").Value = -6.287176378676221
$ws.Range("AF2").Value = 5.132400194802501
$ws.Range("AG2").Value = 6.287176378676221
$ws.Range("AH2").Value = -11.41957657347872
$ws.Range("AI2").Value = 5.132400194802501
$ws.Range("AJ2").Value = 6.287176378676221
$ws.Range("AK2").Value = -5.132400194802501
$ws.Range("AL2").Value = 12.32744885455886
$ws.Range("AM2").Value = -8.214912606512518
$ws.Range("AN2").Value = -45.53134160081841
$ws.Range("AO2").Value = -70.18623930993614
$ws.Range("AP2").Value = -53.75641409691111
$ws.Range("AQ2").Value = -12.32744885455886
$ws.Range("AR2").Value = 8.214912606512518
$ws.Range("AS2").Value = 60.92981298656036
$ws.Range("AT2").Value = 60.92981298656036
$ws.Range("AU2").Value = 61.37293967784092
$ws.Range("AV2").Value = 61.37293967784092
$ws.Range("AW2").Value = 61.37293967784098
$ws.Range("AX2").Value = 60.83066004025909
$ws.Range("AY2").Value = 60.83066004025909
$ws.Range("AZ2").Value = 45.0390723501958
$ws.Range("BA2").Value = 45.0390723501958
$ws.Range("BB2").Value = 44.69642864647369
$ws.Range("BC2").Value = 44.74254362597861
$ws.Range("BD2").Value = 44.63993792111262
$ws.Range("BE2").Value = 45.00241456369662
$ws.Range("BF2").Value = 45.00241456369662
$ws.Range("BG2").Value = 54.99078891686219
$ws.Range("BJ2").Value = 45.95613052988386
$ws.Range("BK2").Value = 53.12368901672909
